$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AR (44) holds the "valor" label for each row; update the value
# from 25 to 18 for every data row (rows 2-51).
$ws.Range("AR2:AR51").Value = 18
